$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 241; this shifts the existing rows 241-295 down to 242-296,
# carrying all of their data (and the dimension) along with them.
$ws.Range("A241:R241").EntireRow.Insert()

# Populate the newly-inserted row 241 with the new record.
$ws.Range("A241").Value = 3
$ws.Range("B241").Value = "Femacal de La Calera"
$ws.Range("C241").Value = "Coquimbo"
$ws.Range("D241").Value = 44711
$ws.Range("E241").Value = 5
$ws.Range("F241").Value = 100112001
$ws.Range("G241").Value = "Berenjena"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 115
$ws.Range("K241").Value = 6500
$ws.Range("L241").Value = 7000
$ws.Range("M241").Value = 6761
$ws.Range("N241").Value = "$/caja 60 unidades"
$ws.Range("O241").Value = "Región de Arica y Parinacota"
$ws.Range("P241").Value = 113
$ws.Range("Q241").Value = 60
$ws.Range("R241").Value = "Hortaliza"
